$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update path-related string cells to absolute Windows paths
$ws.Range("B6").Value = "C:\Users\naeim\OneDrive\Desktop\REWET\Examples\Net3\Result"
$ws.Range("B7").Value = "C:\Users\naeim\OneDrive\Desktop\REWET\Examples\Net3\RunFiles"
$ws.Range("B15").Value = "C:\Users\naeim\OneDrive\Desktop\REWET\Examples\Net3\Net3.inp"
$ws.Range("B19").Value = "C:\Users\naeim\OneDrive\Desktop\REWET\test\test_data\10_day_Net3_No_restoration\test_list.xlsx"
$ws.Range("B20").Value = "C:\Users\naeim\OneDrive\Desktop\REWET\test\test_data\10_day_Net3_No_restoration\Damages"
$ws.Range("B48").Value = "C:\Users\naeim\OneDrive\Desktop\REWET\Examples\Net3\config.txt"

# Flip boolean flags from TRUE to FALSE
$ws.Range("B9").Value = $false
$ws.Range("B23").Value = $false
